$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Apply formatting for the two new columns (N, O) by copying the
#     existing styles from column M (header style for row 1, data style for rows 2-77) ---
$ws.Range("M1").Copy()
$ws.Range("N1:O1").PasteSpecial(-4122)
$ws.Range("M2").Copy()
$ws.Range("N2:O77").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Set the new column widths (N=14, O=15) to match column C (~13.84 chars) ---
$ws.Range("N1:O77").ColumnWidth = 13.14

# --- 3. Header row (row 1): new column titles ---
$ws.Cells.Item(1, 14).Value = "Bezeichnung für Statusübersicht De"
$ws.Cells.Item(1, 15).Value = "Bezeichnung für Statusübersicht En"

# --- 4. Small text corrections elsewhere in the sheet ---
$ws.Cells.Item(21, 4).Value = "b) Wirtschaft"
$ws.Cells.Item(21, 5).Value = "b) Business"
$ws.Cells.Item(73, 13).Value = "Corruption Perception Index in Germany"

# --- 5. Fill the new N/O columns for data rows 2-77 ---
$ws.Cells.Item(2, 14).Value = ""
$ws.Cells.Item(2, 15).Value = ""
$ws.Cells.Item(3, 14).Value = ""
$ws.Cells.Item(3, 15).Value = ""
$ws.Cells.Item(4, 14).Value = ""
$ws.Cells.Item(4, 15).Value = ""
$ws.Cells.Item(5, 14).Value = ""
$ws.Cells.Item(5, 15).Value = ""
$ws.Cells.Item(6, 14).Value = ""
$ws.Cells.Item(6, 15).Value = ""
$ws.Cells.Item(7, 14).Value = ""
$ws.Cells.Item(7, 15).Value = ""
$ws.Cells.Item(8, 14).Value = ""
$ws.Cells.Item(8, 15).Value = ""
$ws.Cells.Item(9, 14).Value = ""
$ws.Cells.Item(9, 15).Value = ""
$ws.Cells.Item(10, 14).Value = ""
$ws.Cells.Item(10, 15).Value = ""
$ws.Cells.Item(11, 14).Value = ""
$ws.Cells.Item(11, 15).Value = ""
$ws.Cells.Item(12, 14).Value = ""
$ws.Cells.Item(12, 15).Value = ""
$ws.Cells.Item(13, 14).Value = ""
$ws.Cells.Item(13, 15).Value = ""
$ws.Cells.Item(14, 14).Value = ""
$ws.Cells.Item(14, 15).Value = ""
$ws.Cells.Item(15, 14).Value = ""
$ws.Cells.Item(15, 15).Value = ""
$ws.Cells.Item(16, 14).Value = ""
$ws.Cells.Item(16, 15).Value = ""
$ws.Cells.Item(17, 14).Value = ""
$ws.Cells.Item(17, 15).Value = ""
$ws.Cells.Item(18, 14).Value = ""
$ws.Cells.Item(18, 15).Value = ""
$ws.Cells.Item(19, 14).Value = ""
$ws.Cells.Item(19, 15).Value = ""
$ws.Cells.Item(20, 14).Value = ""
$ws.Cells.Item(20, 15).Value = ""
$ws.Cells.Item(21, 14).Value = "Frauen in Führungspositionen in der Wirtschaft"
$ws.Cells.Item(21, 15).Value = "Women in management positions in business"
$ws.Cells.Item(22, 14).Value = "Frauen in Führungspositionen im öffentlichen Dienst des Bundes"
$ws.Cells.Item(22, 15).Value = "Women in management positions in the federal civil service"
$ws.Cells.Item(23, 14).Value = ""
$ws.Cells.Item(23, 15).Value = ""
$ws.Cells.Item(24, 14).Value = ""
$ws.Cells.Item(24, 15).Value = ""
$ws.Cells.Item(25, 14).Value = ""
$ws.Cells.Item(25, 15).Value = ""
$ws.Cells.Item(26, 14).Value = ""
$ws.Cells.Item(26, 15).Value = ""
$ws.Cells.Item(27, 14).Value = ""
$ws.Cells.Item(27, 15).Value = ""
$ws.Cells.Item(28, 14).Value = ""
$ws.Cells.Item(28, 15).Value = ""
$ws.Cells.Item(29, 14).Value = ""
$ws.Cells.Item(29, 15).Value = ""
$ws.Cells.Item(30, 14).Value = ""
$ws.Cells.Item(30, 15).Value = ""
$ws.Cells.Item(31, 14).Value = ""
$ws.Cells.Item(31, 15).Value = ""
$ws.Cells.Item(32, 14).Value = ""
$ws.Cells.Item(32, 15).Value = ""
$ws.Cells.Item(33, 14).Value = ""
$ws.Cells.Item(33, 15).Value = ""
$ws.Cells.Item(34, 14).Value = ""
$ws.Cells.Item(34, 15).Value = ""
$ws.Cells.Item(35, 14).Value = ""
$ws.Cells.Item(35, 15).Value = ""
$ws.Cells.Item(36, 14).Value = ""
$ws.Cells.Item(36, 15).Value = ""
$ws.Cells.Item(37, 14).Value = ""
$ws.Cells.Item(37, 15).Value = ""
$ws.Cells.Item(38, 14).Value = ""
$ws.Cells.Item(38, 15).Value = ""
$ws.Cells.Item(39, 14).Value = ""
$ws.Cells.Item(39, 15).Value = ""
$ws.Cells.Item(40, 14).Value = ""
$ws.Cells.Item(40, 15).Value = ""
$ws.Cells.Item(41, 14).Value = ""
$ws.Cells.Item(41, 15).Value = ""
$ws.Cells.Item(42, 14).Value = ""
$ws.Cells.Item(42, 15).Value = ""
$ws.Cells.Item(43, 14).Value = ""
$ws.Cells.Item(43, 15).Value = ""
$ws.Cells.Item(44, 14).Value = ""
$ws.Cells.Item(44, 15).Value = ""
$ws.Cells.Item(45, 14).Value = ""
$ws.Cells.Item(45, 15).Value = ""
$ws.Cells.Item(46, 14).Value = ""
$ws.Cells.Item(46, 15).Value = ""
$ws.Cells.Item(47, 14).Value = ""
$ws.Cells.Item(47, 15).Value = ""
$ws.Cells.Item(48, 14).Value = ""
$ws.Cells.Item(48, 15).Value = ""
$ws.Cells.Item(49, 14).Value = ""
$ws.Cells.Item(49, 15).Value = ""
$ws.Cells.Item(50, 14).Value = ""
$ws.Cells.Item(50, 15).Value = ""
$ws.Cells.Item(51, 14).Value = ""
$ws.Cells.Item(51, 15).Value = ""
$ws.Cells.Item(52, 14).Value = ""
$ws.Cells.Item(52, 15).Value = ""
$ws.Cells.Item(53, 14).Value = ""
$ws.Cells.Item(53, 15).Value = ""
$ws.Cells.Item(54, 14).Value = ""
$ws.Cells.Item(54, 15).Value = ""
$ws.Cells.Item(55, 14).Value = ""
$ws.Cells.Item(55, 15).Value = ""
$ws.Cells.Item(56, 14).Value = ""
$ws.Cells.Item(56, 15).Value = ""
$ws.Cells.Item(57, 14).Value = ""
$ws.Cells.Item(57, 15).Value = ""
$ws.Cells.Item(58, 14).Value = ""
$ws.Cells.Item(58, 15).Value = ""
$ws.Cells.Item(59, 14).Value = ""
$ws.Cells.Item(59, 15).Value = ""
$ws.Cells.Item(60, 14).Value = ""
$ws.Cells.Item(60, 15).Value = ""
$ws.Cells.Item(61, 14).Value = ""
$ws.Cells.Item(61, 15).Value = ""
$ws.Cells.Item(62, 14).Value = ""
$ws.Cells.Item(62, 15).Value = ""
$ws.Cells.Item(63, 14).Value = ""
$ws.Cells.Item(63, 15).Value = ""
$ws.Cells.Item(64, 14).Value = "Stickstoffeintrag über die Zuflüsse in die Ostsee"
$ws.Cells.Item(64, 15).Value = "Nitrogen inputs via the inflows into the Baltic Sea"
$ws.Cells.Item(65, 14).Value = "Stickstoffeintrag über die Zuflüsse in die Nordsee"
$ws.Cells.Item(65, 15).Value = "Nitrogen inputs via the inflows into the North Sea"
$ws.Cells.Item(66, 14).Value = ""
$ws.Cells.Item(66, 15).Value = ""
$ws.Cells.Item(67, 14).Value = ""
$ws.Cells.Item(67, 15).Value = ""
$ws.Cells.Item(68, 14).Value = ""
$ws.Cells.Item(68, 15).Value = ""
$ws.Cells.Item(69, 14).Value = ""
$ws.Cells.Item(69, 15).Value = ""
$ws.Cells.Item(70, 14).Value = ""
$ws.Cells.Item(70, 15).Value = ""
$ws.Cells.Item(71, 14).Value = ""
$ws.Cells.Item(71, 15).Value = ""
$ws.Cells.Item(72, 14).Value = ""
$ws.Cells.Item(72, 15).Value = ""
$ws.Cells.Item(73, 14).Value = ""
$ws.Cells.Item(73, 15).Value = ""
$ws.Cells.Item(74, 14).Value = "Corruption Perception Index in Partnerländern deutscher Entwicklungszusammenarbeit"
$ws.Cells.Item(74, 15).Value = "Corruption Perception Index in partner countries involved in German development cooperation"
$ws.Cells.Item(75, 14).Value = ""
$ws.Cells.Item(75, 15).Value = ""
$ws.Cells.Item(76, 14).Value = ""
$ws.Cells.Item(76, 15).Value = ""
$ws.Cells.Item(77, 14).Value = ""
$ws.Cells.Item(77, 15).Value = ""
